# "Redoing all the effects." - rebuild the Card Effect lookup table on sheet1.
# The table is rewritten column-by-column (A..J) so brand new lookup values
# (new Action Target entries in F, the whole new Colors list in J) land in
# the shared-string table in the same order the author's edit produced them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: Trigger
$colA = @("Trigger", "[SPONTANEOUS]", "[RETALIATE]", "[MANUAL]")
for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $colA[$i]
}

# Column B: Spontaneous Trigger
$colB = @("Spontaneous Trigger", "(when played)", "(on ambush)", "(when drawn)", "(when milled)", "(when vanished)", "(when disintegrated)", "(when restored)", "(when discarded)")
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $colB[$i]
}

# Column C: Phase (row 2 stays blank)
$ws.Cells.Item(1, 3).Value = "Phase"
$colC = @("[PREPARE PHASE]", "[PLAN PHASE]", "[ESCAPE PHASE]", "[RESOLUTION PHASE]")
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item($i + 3, 3).Value = $colC[$i]
}

# Column D: Restriction
$colD = @("Restriction", "[RESTRICT 1]", "[RESTRICT 2]", "[RESTRICT 3]")
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Cells.Item($i + 1, 4).Value = $colD[$i]
}

# Column E: Action
$colE = @("Action", "Draw", "Vanish", "Disintegrate", "Mill", "Restore")
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item($i + 1, 5).Value = $colE[$i]
}

# Column F: Action Target (the long list - cards/anomalies/rooms/activators/actions per color)
$colF = @(
    "Action Target",
    "1 card",
    "2 cards",
    "1 anomaly",
    "2 anomalies",
    "1 room",
    "2 rooms",
    "1 activator",
    "1 action",
    "2 actions",
    "1 <specific card tag> card",
    "1 <specific card tag> card and 1 <specific card tag> card",
    "1 <specific card tag> card or 1 <specific card tag> card",
    "1 red card",
    "1 red anomaly",
    "1 red room",
    "1 red activator",
    "1 red action",
    "1 green card",
    "1 green anomaly",
    "1 green room",
    "1 green activator",
    "1 green action",
    "1 orange card",
    "1 orange anomaly",
    "1 orange room",
    "1 orange activator",
    "1 orange action",
    "1 purple card",
    "1 purple anomaly",
    "1 purple room",
    "1 purple activator",
    "1 purple action",
    "1 gray card",
    "1 gray anomaly",
    "1 gray room",
    "1 gray activator",
    "1 gray action",
    "1 black card",
    "1 black anomaly",
    "1 black room",
    "1 black activator",
    "1 black action"
)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item($i + 1, 6).Value = $colF[$i]
}

# Column G: Downside (row 2 stays blank)
$ws.Cells.Item(1, 7).Value = "Downside"
$colG = @("discard 1 card", "exhaust 1 card you own", "exhaust 2 cards you own", "exhaust 1 anomaly you own", "exhaust 1 activator you own", "vanish 1 card you own", "disintegrate 1 card you own")
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws.Cells.Item($i + 3, 7).Value = $colG[$i]
}

# Column H: Other Keywords (row 2 stays blank)
$ws.Cells.Item(1, 8).Value = "Other Keywords"
$ws.Cells.Item(3, 8).Value = "[DOUBLE AMBUSH]"
$ws.Cells.Item(4, 8).Value = "[IMMUNE]"

# Column I: Card Type
$colI = @("Card Type", "Anomaly", "Room", "Activator", "Action")
for ($i = 0; $i -lt $colI.Length; $i++) {
    $ws.Cells.Item($i + 1, 9).Value = $colI[$i]
}

# Column J: Colors (brand new column)
$colJ = @("Colors", "Red", "Green", "Orange", "Purple", "Gray", "Black")
for ($i = 0; $i -lt $colJ.Length; $i++) {
    $ws.Cells.Item($i + 1, 10).Value = $colJ[$i]
}

# Selection moves to J4 (the view's prior scroll-down to row 4 is dropped
# along with it, since the rebuilt table is short enough to view from A1).
$ws.Range("J4").Select()
